# "Ajout des sprites animes" - add a second "Calendrier" block (Version 6) below
# the existing planning table, plus a "Calendrier Version 5" title above it, by
# duplicating the table header/criteria rows and filling in the controller/
# animated-sprites related content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New title rows -------------------------------------------------------
$ws.Range("A6").Value = "Calendrier Version 5"
$ws.Range("A6").Font.Bold = $true
$ws.Range("A6").VerticalAlignment = -4160

$ws.Range("A20").Value = "Calendrier Version 6"
$ws.Range("A20").Font.Bold = $true
$ws.Range("A20").VerticalAlignment = -4160

# --- Duplicate the "Calendrier" table (rows 9-12) down to rows 22-25 ------
$ws.Range("A9:K9").Copy()
$ws.Range("A22:K22").PasteSpecial(-4122)

$ws.Range("A10:K10").Copy()
$ws.Range("A23:K23").PasteSpecial(-4122)

$ws.Range("A11:K11").Copy()
$ws.Range("A24:K24").PasteSpecial(-4122)

$ws.Range("A12:K12").Copy()
$ws.Range("A25:K25").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Row heights for the (re)wrapped rows
$ws.Rows.Item(9).RowHeight = 57.6
$ws.Rows.Item(10).RowHeight = 86.4
$ws.Rows.Item(15).RowHeight = 28.8
$ws.Rows.Item(22).RowHeight = 43.2
$ws.Rows.Item(23).RowHeight = 72

# --- Fill in the new header row (22): controller / animated sprites -------
$ws.Cells.Item(22, 2).Value = "Sprites animés"
$ws.Cells.Item(22, 3).Value = "Gestion de la manette"
$ws.Cells.Item(22, 4).Value = "Son minimal"
$ws.Cells.Item(22, 5).Value = "Caméra qui suit les joueurs"
$ws.Cells.Item(22, 6).Value = "Combos joueur"
$ws.Cells.Item(22, 7).Value = "Modifier le moteur pour améliorer le gameplay"
$ws.Cells.Item(22, 8).Value = "Musique de fond"
$ws.Cells.Item(22, 9).Value = "S'accrocher a des plates-formes"

# --- Fill in the new "Critères de succès" row (23) -------------------------
$ws.Cells.Item(23, 2).Value = "Animations du personnage en idle`nAnimation d'une attaque spéciale"
$ws.Cells.Item(23, 3).Value = "Pour bouger le personnage avec une manette"
$ws.Cells.Item(23, 4).Value = "1 effet quand fait une attaque, 1 quand touché"
$ws.Cells.Item(23, 5).Value = "Joueurs toujours à l'écran, Effet de zoom arrière limité"
$ws.Cells.Item(23, 6).Value = "1 attaque spéciale avec un combo particulier (exemple, droite+A appuyé pendant 2 secondes)"
$ws.Cells.Item(23, 7).Value = "Rebondir sur les murs"
$ws.Cells.Item(23, 8).Value = "Pouvoir choisir entre 2 musiques de fond"

# --- Clear the numeric "Importance"/"Difficulté" values that were copied ---
# from rows 11/12 but should stay blank in the new table (rows 24/25)
$ws.Range("B24:K24").ClearContents()
$ws.Range("B25:K25").ClearContents()

# --- Trim the copied ranges down to the real target shape ------------------
$ws.Cells.Item(22, 10).Clear()   # J22
$ws.Cells.Item(23, 10).Clear()   # J23
$ws.Cells.Item(23, 11).Clear()   # K23
$ws.Cells.Item(24, 3).Clear()    # C24
$ws.Cells.Item(24, 11).Clear()   # K24
$ws.Cells.Item(25, 3).Clear()    # C25
$ws.Cells.Item(25, 11).Clear()   # K25

# --- Plain formatted (blank) rows around/below the new table --------------
$ws.Range("B20:I20").VerticalAlignment = -4160
$ws.Range("A21:I21").VerticalAlignment = -4160
$ws.Range("K21").VerticalAlignment = -4160
$ws.Range("K22").VerticalAlignment = -4160
$ws.Range("A27:I27").VerticalAlignment = -4160

# --- Clean up the old trailing blank rows (17-19) that got absorbed into ---
# the new layout
$ws.Range("A17:I17").Clear()
$ws.Range("A18:I18").Clear()
$ws.Range("C19:J19").Clear()

# --- Restore selection -------------------------------------------------
$ws.Range("G20").Select()
